$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IAM_PERMISSION")

# Update PATH column (G) for rows 8-22: replace leading "/iam/" with "/asgard/"
for ($row = 8; $row -le 22; $row++) {
    $cell = $ws.Cells.Item($row, 7)  # column G
    $val = $cell.Value2
    if ($val -and $val.ToString().StartsWith("/iam/")) {
        $cell.Value2 = $val.ToString().Replace("/iam/", "/asgard/")
    }
}
